$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 33: raw and clean data for July 2nd (2020-07-02)
# Leading apostrophe forces Excel to store the date-like text as a string
# (matching column A's existing "t=s" shared-string cells) instead of
# auto-converting it to a date serial value. ClearFormats() then removes
# the quote-prefix cell formatting so the cell stays unstyled, like the
# other date cells in column A.
$ws.Cells.Item(33, 1).Value = "'2020-07-02"
$ws.Cells.Item(33, 1).ClearFormats()

$ws.Cells.Item(33, 2).Value = 238511
$ws.Cells.Item(33, 3).Value = 295561
$ws.Cells.Item(33, 4).Value = 76423
$ws.Cells.Item(33, 5).Value = 29189
$ws.Cells.Item(33, 6).Value = 30.52
